$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Future Heads" debug/analysis table in columns N (labels) and O
# (distances) is no longer produced — the food-chasing logic now acts
# immediately when food is close, so the extra reporting columns are
# cleared out.

# Rows 2-5 and 7 keep their label styling (s="3") but no longer carry a
# value; row 2-5 also had a numeric distance in column O which is removed
# entirely.
$ws.Range("N2:N5").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("O2:O5").ClearContents()

# Rows 8-20 held the full list of "Position(x=.., y=..)" future-head
# entries; that list is gone now.
$ws.Range("N8:N20").ClearContents()

# Reflect the new area of interest: the (now empty) report block.
$ws.Range("N2:P23").Select() | Out-Null
